$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous data block (will be rewritten in the new sorted order)
$ws.Range("A2:J15").ClearContents()
$ws.Range("A2:J15").ClearFormats()

# Write row data (columns A, B, D, E, F, G, H, I, J; column C is set via Hyperlinks.Add below)
# Row 2: Afrobarometer
$ws.Range("A2").Value = 'Afrobarometer'
$ws.Range("B2").Value = 'citizens'

# Row 3: American National Election Studies
$ws.Range("A3").Value = 'American National Election Studies'
$ws.Range("B3").Value = 'citizens'
$ws.Range("E3").Value = 'USA'
$ws.Range("F3").Value = 1948
$ws.Range("H3").Value = 'public archive'
$ws.Range("I3").Value = 'free, online'
$ws.Range("J3").Value = 'cross-sectional, panel'

# Row 4: Americas Barometer
$ws.Range("A4").Value = 'Americas Barometer'
$ws.Range("B4").Value = 'citizens'

# Row 5: Asian Barometer
$ws.Range("A5").Value = 'Asian Barometer'
$ws.Range("B5").Value = 'citizens'

# Row 6: Caucasus Barometer
$ws.Range("A6").Value = 'Caucasus Barometer'
$ws.Range("B6").Value = 'citizens'

# Row 7: Comparative Political Data Set
$ws.Range("A7").Value = 'Comparative Political Data Set'
$ws.Range("B7").Value = 'political institutions'

# Row 8: Comparative Study of Electoral Systems
$ws.Range("A8").Value = 'Comparative Study of Electoral Systems'
$ws.Range("B8").Value = 'citizens'

# Row 9: Correlates of War
$ws.Range("A9").Value = 'Correlates of War'
$ws.Range("B9").Value = 'international relations'
$ws.Range("E9").Value = 'world'
$ws.Range("F9").Value = 1816
$ws.Range("G9").Value = 2007
$ws.Range("H9").Value = 'public archive'
$ws.Range("I9").Value = 'free, online'

# Row 10: Democratic Accountability and Citizen-Politician Linkages
$ws.Range("A10").Value = 'Democratic Accountability and Citizen-Politician Linkages'
$ws.Range("B10").Value = 'democracy'

# Row 11: Eurobarometer
$ws.Range("A11").Value = 'Eurobarometer'
$ws.Range("B11").Value = 'citizens'

# Row 12: European Quality of Life
$ws.Range("A12").Value = 'European Quality of Life'
$ws.Range("B12").Value = 'citizens'

# Row 13: European Social Survey
$ws.Range("A13").Value = 'European Social Survey'
$ws.Range("B13").Value = 'citizens'
$ws.Range("D13").Value = 'Media and social trust; Subjective well-being; Human values; Immigration; Citizen involvement; Health and care; Economic morality; Family work and well-being; Timing of life; Personal well-being; Welfare attitudes; Ageism; Justice; Democracy; Social inequalities in health; Public attitudes to climate change'
$ws.Range("E13").Value = 'Europe'
$ws.Range("F13").Value = 2002
$ws.Range("H13").Value = 'public archive'
$ws.Range("I13").Value = 'free, online'
$ws.Range("J13").Value = 'cross-sectional'

# Row 14: European Values Study
$ws.Range("A14").Value = 'European Values Study'
$ws.Range("B14").Value = 'citizens'

# Row 15: International Social Survey Programme
$ws.Range("A15").Value = 'International Social Survey Programme'
$ws.Range("B15").Value = 'citizens'

# Row 16: Longitudinal Internet Studies for the Social sciences
$ws.Range("A16").Value = 'Longitudinal Internet Studies for the Social sciences'
$ws.Range("B16").Value = 'citizens'
$ws.Range("D16").Value = 'Health; Religion and Ethnicity; Social Integration and Leisure; Family and Household; Work and Schooling; Personality; Politics and Values; Economic Situation'
$ws.Range("E16").Value = 'Netherlands'
$ws.Range("F16").Value = 2007
$ws.Range("H16").Value = 'public archive'
$ws.Range("I16").Value = 'free, letter'
$ws.Range("J16").Value = 'panel'

# Row 17: Parties, Governments and Legislatures Dataset
$ws.Range("A17").Value = 'Parties, Governments and Legislatures Dataset'
$ws.Range("B17").Value = 'governments and parties'

# Row 18: Pew Global Attitudes Survey
$ws.Range("A18").Value = 'Pew Global Attitudes Survey'
$ws.Range("B18").Value = 'citizens'
$ws.Range("D18").Value = 'Attitudes'
$ws.Range("E18").Value = 'world'
$ws.Range("F18").Value = 2001
$ws.Range("H18").Value = 'public archive'
$ws.Range("I18").Value = 'free, online'
$ws.Range("J18").Value = 'cross-sectional'

# Row 19: World Values Survey
$ws.Range("A19").Value = 'World Values Survey'
$ws.Range("B19").Value = 'citizens'
$ws.Range("D19").Value = 'Political attitudes'
$ws.Range("E19").Value = 'world'
$ws.Range("F19").Value = 1981
$ws.Range("H19").Value = 'public archive'
$ws.Range("I19").Value = 'free, online'
$ws.Range("J19").Value = 'cross-sectional'

# Column C: hyperlinks (also sets the display text + hyperlink style)
$ws.Hyperlinks.Add($ws.Range("C2"), 'http://www.afrobarometer.org/', "", "", 'http://www.afrobarometer.org/')
$ws.Hyperlinks.Add($ws.Range("C3"), 'http://www.electionstudies.org/', "", "", 'http://www.electionstudies.org/')
$ws.Hyperlinks.Add($ws.Range("C4"), 'http://www.vanderbilt.edu/lapop/about-americasbarometer.php', "", "", 'http://www.vanderbilt.edu/lapop/about-americasbarometer.php')
$ws.Hyperlinks.Add($ws.Range("C5"), 'http://www.asianbarometer.org/', "", "", 'http://www.asianbarometer.org/')
$ws.Hyperlinks.Add($ws.Range("C6"), 'http://caucasusbarometer.org/en/datasets/', "", "", 'http://caucasusbarometer.org/en/datasets/')
$ws.Hyperlinks.Add($ws.Range("C7"), 'http://www.cpds-data.org/', "", "", 'http://www.cpds-data.org/')
$ws.Hyperlinks.Add($ws.Range("C8"), 'http://www.cses.org/', "", "", 'http://www.cses.org/')
$ws.Hyperlinks.Add($ws.Range("C9"), 'http://www.correlatesofwar.org/', "", "", 'http://www.correlatesofwar.org/')
$ws.Hyperlinks.Add($ws.Range("C10"), 'https://sites.duke.edu/democracylinkage/data/', "", "", 'https://sites.duke.edu/democracylinkage/data/')
$ws.Hyperlinks.Add($ws.Range("C11"), 'http://ec.europa.eu/commfrontoffice/publicopinion/index.cfm', "", "", 'http://ec.europa.eu/commfrontoffice/publicopinion/index.cfm')
$ws.Hyperlinks.Add($ws.Range("C12"), 'https://www.eurofound.europa.eu/surveys/european-quality-of-life-surveys', "", "", 'https://www.eurofound.europa.eu/surveys/european-quality-of-life-surveys')
$ws.Hyperlinks.Add($ws.Range("C13"), 'http://www.europeansocialsurvey.org/', "", "", 'http://www.europeansocialsurvey.org/')
$ws.Hyperlinks.Add($ws.Range("C14"), 'http://www.europeanvaluesstudy.eu/', "", "", 'http://www.europeanvaluesstudy.eu/')
$ws.Hyperlinks.Add($ws.Range("C15"), 'http://www.issp.org/menu-top/home/', "", "", 'http://www.issp.org/menu-top/home/')
$ws.Hyperlinks.Add($ws.Range("C16"), 'https://www.lissdata.nl/', "", "", 'https://www.lissdata.nl/')
$ws.Hyperlinks.Add($ws.Range("C17"), 'http://www.edac.eu/policies_desc.cfm?v_id=112', "", "", 'http://www.edac.eu/policies_desc.cfm?v_id=112')
$ws.Hyperlinks.Add($ws.Range("C18"), 'http://www.pewglobal.org/datasets/', "", "", 'http://www.pewglobal.org/datasets/')
$ws.Hyperlinks.Add($ws.Range("C19"), 'http://www.worldvaluessurvey.org/wvs.jsp', "", "", 'http://www.worldvaluessurvey.org/wvs.jsp')

# Re-apply vertical-center alignment style to the rows that carried it before, at their new positions
$ws.Range("A8").VerticalAlignment = -4108
$ws.Range("A14").VerticalAlignment = -4108
$ws.Range("A15").VerticalAlignment = -4108
$ws.Range("D19").VerticalAlignment = -4108

# Update selection to match the final state
$ws.Range("C23").Select()
